# Applies the "Updated cryptos list" data refresh to Sheet1.
# Each data row (2-51) holds Coin / Link / Price / Volume(1h) as plain text
# (inline/shared strings), so values are written back as text too. For cells
# whose new text happens to look like a plain number (e.g. "304.36"), the
# NumberFormat is temporarily switched to Text ("@") before the assignment so
# Excel does not silently reinterpret the string as a numeric value; the
# formatting is cleared again right after so the cell keeps its original
# (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.853.32"
$ws.Range("E2").Value = "  +0.39%  "
# Row 3
$ws.Range("D3").Value = "2.548.15"
$ws.Range("E3").Value = "  -0.12%  "
# Row 4
$ws.Range("E4").Value = "  -0.11%  "
# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "304.36"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.50%  "
# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "97.90"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +5.86%  "
# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.578"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.32%  "
# Row 8
$ws.Range("E8").Value = "  +0.11%  "
# Row 9
$ws.Range("E9").Value = "  -0.58%  "
# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.88"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +2.33%  "
# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0825"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +2.62%  "
# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.94%  "
# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.59"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.79%  "
# Row 14
$ws.Range("D14").Value = "2.939.32"
$ws.Range("E14").Value = "  -0.04%  "
# Row 15
$ws.Range("D15").Value = "2.538.15"
$ws.Range("E15").Value = "  +0.68%  "
# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "15.09"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +5.26%  "
# Row 17
$ws.Range("E17").Value = "  -1.15%  "
# Row 18
$ws.Range("D18").Value = "42.864.72"
$ws.Range("E18").Value = "  +0.37%  "
# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.33"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +2.95%  "
# Row 20
$ws.Range("E20").Value = "  +0.34%  "
# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -0.58%  "
# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "71.83"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.04%  "
# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "254.49"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.92%  "
# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.95"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.95%  "
# Row 25
$ws.Range("E25").Value = "  -2.25%  "
# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "28.01"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -4.07%  "
# Row 27
$ws.Range("E27").Value = "  +0.10%  "
# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.28"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +7.60%  "
# Row 29
$ws.Range("E29").Value = "  +0.93%  "
# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "37.98"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +2.17%  "
# Row 31
$ws.Range("E31").Value = "  +2.34%  "
# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "156.97"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +2.99%  "
# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "19.43"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +14.78%  "
# Row 34
$ws.Range("E34").Value = "  -0.95%  "
# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.33"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -2.10%  "
# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0798"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +0.51%  "
# Row 37
$ws.Range("E37").Value = "  -4.64%  "
# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.115"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.89%  "
# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "24.65"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.11%  "
# Row 40
$ws.Range("E40").Value = "  +0.54%  "
# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +31.36%  "
# Row 42
$ws.Range("E42").Value = "  +0.42%  "
# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.87"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.23%  "
# Row 44
$ws.Range("D44").Value = "2.100.43"
$ws.Range("E44").Value = "  +0.77%  "
# Row 45
$ws.Range("E45").Value = "  -1.99%  "
# Row 46
$ws.Range("E46").Value = "  +0.07%  "
# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "87.07"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +2.37%  "
# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.89"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -1.67%  "
# Row 49
$ws.Range("D49").Value = "2.796.77"
$ws.Range("E49").Value = "  +0.04%  "
# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "74.26"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +7.33%  "
# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.192"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +2.05%  "
